# Update DeudoresPrueba worksheet: refresh client/date/amount data for rows 4-34
# (rows 2-3 are unchanged), then drop the two trailing rows (35-36) that no
# longer exist in the refreshed data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consecutivo (A), Cliente (B), Fecha serial (C), Valor (D) for each data row.
# Row 2 and 3 keep their original values (ALISO / ARROZ PAISA SUBA) untouched.
$data = @(
    @(4,  3,  "CAMPO VERDE ZIPAQUIRA", 46014, 71900),
    @(5,  4,  "CANTON WOK",            46015, 252000),
    @(6,  5,  "CARNES JOHANA",         46011, 166000),
    @(7,  6,  "CIMARRON DORADO",       46010, 375000),
    @(8,  7,  "CIMARRON DORADO",       46017, 315600),
    @(9,  8,  "CLIENTE PAOLA",         46018, 274000),
    @(10, 9,  "COCINA CHINA",          46018, 170000),
    @(11, 10, "CRISTIAN ACACIAS",      46009, 1000000),
    @(12, 11, "DARWIN FUTBOL",         45921, 200000),
    @(13, 12, "DAVIDCITO",             45947, 100000),
    @(14, 13, "FRANCO",                46017, 545800),
    @(15, 14, "FRANCO",                45996, 20000),
    @(16, 15, "LA PAMPA",              46006, 229900),
    @(17, 16, "LA SELECTA",            45912, 82000),
    @(18, 17, "MAFE",                  46017, 190000),
    @(19, 18, "MERKA FRUVER DEXI",     45988, 15400),
    @(20, 19, "MERKA FRUVER DEXI",     45995, 339000),
    @(21, 20, "MICHAEL",               46011, 80000),
    @(22, 21, "NEVADA",                46017, 195000),
    @(23, 22, "NOVILLON SAN MATEO",    45971, 33000),
    @(24, 23, "PARAÍSO MOSQUERA",      46013, 328800),
    @(25, 24, "PINILLA",               45931, 82000),
    @(26, 25, "PINILLA SOACHA",        46015, 166000),
    @(27, 26, "PLAZA JESSICA",         46014, 1655400),
    @(28, 27, "PUNTA DE ANCA",         46017, 507000),
    @(29, 28, "SAMY 2",                46013, 142000),
    @(30, 29, "SAN JOAQUIN",           46015, 229300),
    @(31, 30, "SANTANDER SUR",         46014, 253000),
    @(32, 31, "SANTANDER SUR",         46018, 218000),
    @(33, 32, "TIMO",                  46015, 189000),
    @(34, 33, "WILINTONG",             46006, 150000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# The refreshed data set only has 33 rows (2-34); remove the now-obsolete
# rows 35 and 36 from the bottom of the sheet.
$ws.Rows("35:36").Delete()
